$d = $word.ActiveDocument
$wns = 'xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"'

# ---------------------------------------------------------------------------
# Hunk 1: merge the four runs of the "Séptima Reunión (19/04/2022) " heading
# into a single run (paragraph 34). A simple whole-paragraph Find/Replace is
# fine here because ALL runs in that paragraph should end up merged anyway.
# ---------------------------------------------------------------------------
$pHeading7 = $d.Paragraphs.Item(34)
$pHeading7.Range.Find.Execute("Séptima Reunión (19/04/2022) ", $false, $false, $false, $false, $false, $true, 1, $false, "Séptima Reunión (19/04/2022) ", 2) | Out-Null

# ---------------------------------------------------------------------------
# Hunk 2: merge "Puesta en común ... la " + "fase 2" into a single run, while
# leaving the following "." run untouched (paragraph 36).
#
# Any in-place edit (Find/Replace, Delete, Text assignment) touching a
# paragraph causes the engine to renormalize/merge ALL sibling runs that
# share formatting - including the "." run we want to keep separate. To
# avoid that, we delete the paragraph's text (keeping its paragraph mark/
# pPr) and then insert a fresh two-run paragraph via InsertXML, which
# preserves the exact run layout we specify.
# ---------------------------------------------------------------------------
$pFase2 = $d.Paragraphs.Item(36)
$s2 = $pFase2.Range.Start
$e2 = $pFase2.Range.End
$d.Range($s2, $e2 - 1).Delete()
$insert2 = $d.Range($s2, $s2)
$xml2 = '<w:p ' + $wns + '>' + `
  '<w:r><w:t>Puesta en común de todos los avances desarrollados a lo largo de la fase 2</w:t></w:r>' + `
  '<w:r><w:t>.</w:t></w:r>' + `
  '</w:p>'
$insert2.InsertXML($xml2)

Write-Host "Hunks 1-2 applied. Paragraph count: $($d.Paragraphs.Count)"

# ---------------------------------------------------------------------------
# Hunk 3: Resumen de la octava reunión.
#
# Paragraph 38 currently reads (across 3 runs):
#   "Inicio de la fase " + "3" + " de la planificación de desarrollo del proyecto."
# It must become a single run with the same text (simple merge is fine, no
# run needs to stay separate), and new paragraphs describing the eighth
# meeting must be inserted right after it (before the old trailing empty
# paragraph).
# ---------------------------------------------------------------------------
$pFase3 = $d.Paragraphs.Item(38)
$pFase3.Range.Find.Execute("Inicio de la fase 3 de la planificación de desarrollo del proyecto.", $false, $false, $false, $false, $false, $true, 1, $false, "Inicio de la fase 3 de la planificación de desarrollo del proyecto.", 2) | Out-Null

Write-Host "Paragraph 38 merged: [$($d.Paragraphs.Item(38).Range.Text)]"

# Insertion point: right after paragraph 38 (just before its paragraph mark end).
$p38 = $d.Paragraphs.Item(38)
$insertAt = $d.Range($p38.Range.End - 1, $p38.Range.End - 1)

# New empty (centered) paragraph.
$xmlBlank1 = '<w:p ' + $wns + '><w:pPr><w:jc w:val="center"/></w:pPr></w:p>'

# New "Octava Reunión (19/04/2022) " heading paragraph (2 runs, bold/underline sz32).
$xmlHeading8 = '<w:p ' + $wns + '>' + `
  '<w:pPr><w:jc w:val="center"/><w:rPr><w:b/><w:bCs/><w:sz w:val="32"/><w:szCs w:val="32"/><w:u w:val="single"/></w:rPr></w:pPr>' + `
  '<w:r><w:rPr><w:b/><w:bCs/><w:sz w:val="32"/><w:szCs w:val="32"/><w:u w:val="single"/></w:rPr><w:t>Octava</w:t></w:r>' + `
  '<w:r><w:rPr><w:b/><w:bCs/><w:sz w:val="32"/><w:szCs w:val="32"/><w:u w:val="single"/></w:rPr><w:t xml:space="preserve"> Reunión (19/04/2022) </w:t></w:r>' + `
  '</w:p>'

# New "Puesta en común ... fase 3." paragraph (3 runs).
$xmlFase3Puesta = '<w:p ' + $wns + '><w:pPr><w:jc w:val="center"/></w:pPr>' + `
  '<w:r><w:t xml:space="preserve">Puesta en común de todos los avances desarrollados a lo largo de la fase </w:t></w:r>' + `
  '<w:r><w:t>3</w:t></w:r>' + `
  '<w:r><w:t>.</w:t></w:r>' + `
  '</w:p>'

# New "Inicio de la fase 4 ..." paragraph (3 runs).
$xmlFase4Inicio = '<w:p ' + $wns + '><w:pPr><w:jc w:val="center"/></w:pPr>' + `
  '<w:r><w:t xml:space="preserve">Inicio de la fase </w:t></w:r>' + `
  '<w:r><w:t>4</w:t></w:r>' + `
  '<w:r><w:t xml:space="preserve"> de la planificación de desarrollo del proyecto.</w:t></w:r>' + `
  '</w:p>'

# New "[Falta de realización de actividades injustificada] ..." paragraph (5 runs).
$xmlFalta = '<w:p ' + $wns + '><w:pPr><w:jc w:val="center"/></w:pPr>' + `
  '<w:r><w:t xml:space="preserve">[Falta de realización de actividades injustificada]  </w:t></w:r>' + `
  '<w:r><w:t>Daniel Yanel Gorrón</w:t></w:r>' + `
  '<w:r><w:t xml:space="preserve"> y </w:t></w:r>' + `
  '<w:r><w:t>Rubén López</w:t></w:r>' + `
  '<w:r><w:t>.</w:t></w:r>' + `
  '</w:p>'

# New trailing empty (centered) paragraph.
$xmlBlank2 = '<w:p ' + $wns + '><w:pPr><w:jc w:val="center"/></w:pPr></w:p>'

$fullXml = $xmlBlank1 + $xmlHeading8 + $xmlFase3Puesta + $xmlFase4Inicio + $xmlFalta + $xmlBlank2

$insertAt.InsertXML($fullXml)

Write-Host "Final paragraph count: $($d.Paragraphs.Count)"
